$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header: H1 becomes year 2023, I1 becomes the "total" label that used to be in H1
$ws.Range("H1").Value = 2023
$ws.Range("I1").Value = $ws.Range("H1").Value
$ws.Range("I1").Value = "total"

# New 2023 data column (H2:H10)
$ws.Range("H2").Value = 43
$ws.Range("H3").Value = 27
$ws.Range("H4").Value = 12
$ws.Range("H5").Value = 7
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 2
$ws.Range("H8").Value = 0
$ws.Range("H9").Value = 91
$ws.Range("H10").Value = 91

# Totals column moves to I, formulas now sum B:H
$ws.Range("I2").Formula = "=SUM(B2:H2)"
$ws.Range("I3").Formula = "=SUM(B3:H3)"
$ws.Range("I4").Formula = "=SUM(B4:H4)"
$ws.Range("I5").Formula = "=SUM(B5:H5)"
$ws.Range("I6").Formula = "=SUM(B6:H6)"
$ws.Range("I7").Formula = "=SUM(B7:H7)"
$ws.Range("I8").Formula = "=SUM(B8:H8)"
$ws.Range("I9").Formula = "=SUM(B9:H9)"
$ws.Range("I10").Formula = "=SUM(B10:H10)"

# Copy style from the old H column (totals) to the new I column, and apply
# plain numeric style (same as B:G) to the new H data column
$ws.Range("H1:H10").Copy()
$ws.Range("I1:I10").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("B1:G10").Copy()
$ws.Range("H2:H8").PasteSpecial(-4122)
$ws.Range("H10").PasteSpecial(-4122)

# H9/H10 used a slightly different style (s=2) like F9/G9
$ws.Range("G9:G9").Copy()
$ws.Range("H9").PasteSpecial(-4122)

# Re-set values after paste (paste formats shouldn't disturb values but formulas might've been affected)
$ws.Range("H1").Value = 2023
$ws.Range("H2").Value = 43
$ws.Range("H3").Value = 27
$ws.Range("H4").Value = 12
$ws.Range("H5").Value = 7
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 2
$ws.Range("H8").Value = 0
$ws.Range("H9").Value = 91
$ws.Range("H10").Value = 91
$ws.Range("I1").Value = "total"

# Selection as per diff
$ws.Range("I9:I10").Select()

$ws.Calculate()
